$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp (A1)
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 7 de Julio de 2020 a las 11:05"   # A1: 'Datos actualizados a 7 de Julio de 2020 a las 09:48' -> 'Datos actualizados a 7 de Julio de 2020 a las 11:05'

# Row 20: 'Francia' slot now holds 'Banglades' data
$ws.Cells.Item(20, 1).Value = "Banglades"   # A20: 'Francia' -> 'Banglades'
$ws.Cells.Item(20, 2).Value = 168645   # B20: 168335 -> 168645
$ws.Cells.Item(20, 3).Value = 3027   # C20: 0 -> 3027
$ws.Cells.Item(20, 4).Value = 78102   # D20: 77308 -> 78102
$ws.Cells.Item(20, 5).Value = 88392   # E20: 61107 -> 88392
$ws.Cells.Item(20, 7).Value = 55   # G20: 0 -> 55
$ws.Cells.Item(20, 8).Value = 2151   # H20: 29920 -> 2151

# Row 21: 'Banglades' slot now holds 'Francia' data
$ws.Cells.Item(21, 1).Value = "Francia"   # A21: 'Banglades' -> 'Francia'
$ws.Cells.Item(21, 2).Value = 168335   # B21: 165618 -> 168335
$ws.Cells.Item(21, 4).Value = 77308   # D21: 76149 -> 77308
$ws.Cells.Item(21, 5).Value = 61107   # E21: 87373 -> 61107
$ws.Cells.Item(21, 8).Value = 29920   # H21: 2096 -> 29920

# Row 39: 'Oman' slot now holds 'Filipinas' data
$ws.Cells.Item(39, 1).Value = "Filipinas"   # A39: 'Oman' -> 'Filipinas'
$ws.Cells.Item(39, 2).Value = 47873   # B39: 47735 -> 47873
$ws.Cells.Item(39, 3).Value = 1540   # C39: 0 -> 1540
$ws.Cells.Item(39, 4).Value = 12386   # D39: 29146 -> 12386
$ws.Cells.Item(39, 5).Value = 34178   # E39: 18371 -> 34178
$ws.Cells.Item(39, 7).Value = 6   # G39: 0 -> 6
$ws.Cells.Item(39, 8).Value = 1309   # H39: 218 -> 1309

# Row 40: 'Filipinas' slot now holds 'Oman' data
$ws.Cells.Item(40, 1).Value = "Oman"   # A40: 'Filipinas' -> 'Oman'
$ws.Cells.Item(40, 2).Value = 47735   # B40: 46333 -> 47735
$ws.Cells.Item(40, 4).Value = 29146   # D40: 12185 -> 29146
$ws.Cells.Item(40, 5).Value = 18371   # E40: 32845 -> 18371
$ws.Cells.Item(40, 8).Value = 218   # H40: 1303 -> 218

# Row 46: 'Polonia' - updated figures
$ws.Cells.Item(46, 2).Value = 36412   # B46: 36155 -> 36412
$ws.Cells.Item(46, 3).Value = 257   # C46: 0 -> 257
$ws.Cells.Item(46, 4).Value = 24238   # D46: 23966 -> 24238
$ws.Cells.Item(46, 5).Value = 10646   # E46: 10668 -> 10646
$ws.Cells.Item(46, 7).Value = 7   # G46: 0 -> 7
$ws.Cells.Item(46, 8).Value = 1528   # H46: 1521 -> 1528

# Row 49: 'Israel' - updated figures
$ws.Cells.Item(49, 2).Value = 31271   # B49: 31186 -> 31271
$ws.Cells.Item(49, 3).Value = 522   # C49: 437 -> 522
$ws.Cells.Item(49, 5).Value = 12802   # E49: 12717 -> 12802

# Row 60: 'Austria' - updated figures
$ws.Cells.Item(60, 2).Value = 18421   # B60: 18365 -> 18421
$ws.Cells.Item(60, 3).Value = 56   # C60: 0 -> 56
$ws.Cells.Item(60, 4).Value = 16686   # D60: 16647 -> 16686
$ws.Cells.Item(60, 5).Value = 1029   # E60: 1012 -> 1029

# Row 76: 'Kirguistan' slot now holds 'El Salvador' data
$ws.Cells.Item(76, 1).Value = "El Salvador"   # A76: 'Kirguistan' -> 'El Salvador'
$ws.Cells.Item(76, 2).Value = 8307   # B76: 8141 -> 8307
$ws.Cells.Item(76, 3).Value = 280   # C76: 450 -> 280
$ws.Cells.Item(76, 4).Value = 4929   # D76: 2916 -> 4929
$ws.Cells.Item(76, 5).Value = 3149   # E76: 5126 -> 3149
$ws.Cells.Item(76, 7).Value = 6   # G76: 7 -> 6
$ws.Cells.Item(76, 8).Value = 229   # H76: 99 -> 229

# Row 77: 'Kenia' slot now holds 'Kirguistan' data
$ws.Cells.Item(77, 1).Value = "Kirguistan"   # A77: 'Kenia' -> 'Kirguistan'
$ws.Cells.Item(77, 2).Value = 8141   # B77: 8067 -> 8141
$ws.Cells.Item(77, 3).Value = 450   # C77: 0 -> 450
$ws.Cells.Item(77, 4).Value = 2916   # D77: 2414 -> 2916
$ws.Cells.Item(77, 5).Value = 5126   # E77: 5489 -> 5126
$ws.Cells.Item(77, 7).Value = 7   # G77: 0 -> 7
$ws.Cells.Item(77, 8).Value = 99   # H77: 164 -> 99

# Row 78: 'El Salvador' slot now holds 'Kenia' data
$ws.Cells.Item(78, 1).Value = "Kenia"   # A78: 'El Salvador' -> 'Kenia'
$ws.Cells.Item(78, 2).Value = 8067   # B78: 8027 -> 8067
$ws.Cells.Item(78, 4).Value = 2414   # D78: 4785 -> 2414
$ws.Cells.Item(78, 5).Value = 5489   # E78: 3019 -> 5489
$ws.Cells.Item(78, 8).Value = 164   # H78: 223 -> 164

# Row 96: 'Estado de Palestina' - updated figures
$ws.Cells.Item(96, 5).Value = 3832   # E96: 3833 -> 3832
$ws.Cells.Item(96, 7).Value = 1   # G96: 0 -> 1
$ws.Cells.Item(96, 8).Value = 18   # H96: 17 -> 18

# Row 121: 'Eslovenia' - updated figures
$ws.Cells.Item(121, 2).Value = 1739   # B121: 1716 -> 1739
$ws.Cells.Item(121, 3).Value = 23   # C121: 0 -> 23
$ws.Cells.Item(121, 4).Value = 1423   # D121: 1384 -> 1423
$ws.Cells.Item(121, 5).Value = 205   # E121: 221 -> 205

# Row 156: 'Tanzania' slot now holds 'Namibia' data
$ws.Cells.Item(156, 1).Value = "Namibia"   # A156: 'Tanzania' -> 'Namibia'
$ws.Cells.Item(156, 2).Value = 539   # B156: 509 -> 539
$ws.Cells.Item(156, 3).Value = 54   # C156: 0 -> 54
$ws.Cells.Item(156, 4).Value = 25   # D156: 183 -> 25
$ws.Cells.Item(156, 5).Value = 514   # E156: 305 -> 514
$ws.Cells.Item(156, 8).Value = 0   # H156: 21 -> 0

# Row 157: 'Namibia' slot now holds 'Tanzania' data
$ws.Cells.Item(157, 1).Value = "Tanzania"   # A157: 'Namibia' -> 'Tanzania'
$ws.Cells.Item(157, 2).Value = 509   # B157: 485 -> 509
$ws.Cells.Item(157, 4).Value = 183   # D157: 25 -> 183
$ws.Cells.Item(157, 5).Value = 305   # E157: 460 -> 305
$ws.Cells.Item(157, 8).Value = 21   # H157: 0 -> 21

# Row 173: 'Islas Feroe' - updated figures
$ws.Cells.Item(173, 4).Value = 188   # D173: 187 -> 188
$ws.Cells.Item(173, 5).Value = 0   # E173: 1 -> 0

# Row 209: 'Groenlandia' slot now holds 'Islas Malvinas' data
$ws.Cells.Item(209, 1).Value = "Islas Malvinas"   # A209: 'Groenlandia' -> 'Islas Malvinas'

# Row 210: 'Islas Malvinas' slot now holds 'Groenlandia' data
$ws.Cells.Item(210, 1).Value = "Groenlandia"   # A210: 'Islas Malvinas' -> 'Groenlandia'
